$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "ID" row (row 2, "建筑ID") — the building-level Id
# field is no longer part of this struct; remaining rows shift up.
$ws.Rows(2).Delete()

# The Prefab / NormalStateFunc / UpStateFunc / Desc fields are now typed
# as "string" instead of "int" (rows 4-7 after the delete above).
$ws.Range("B4:B7").Value = "string"

# Move the active selection, matching where the author left off editing.
$ws.Range("G14").Select() | Out-Null

# Data validation ranges need to be re-anchored one row up (F9:F1048576 ->
# F8:F1048576) now that row 2 is gone; recreate both rules cleanly.
$ws.Cells.Validation.Delete() | Out-Null
$ws.Range("F1").Validation.Add(0, 1, 1, "") | Out-Null
$ws.Range("F8:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"') | Out-Null
